# Apply updated crypto price/volume figures per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.575.42'
$ws.Range('E2').Value = '  -2.08%  '
$ws.Range('D3').Value = '2.586.24'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.94'
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.95'
$ws.Range('E6').Value = '  -4.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  -5.95%  '
$ws.Range('D9').Value = '2.592.69'
$ws.Range('E9').Value = '  -3.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.60'
$ws.Range('E10').Value = '  +6.14%  '
$ws.Range('E11').Value = '  -3.32%  '
$ws.Range('E12').Value = '  -1.15%  '
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('D14').Value = '3.037.26'
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('D15').Value = '60.590.89'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.72'
$ws.Range('E16').Value = '  -4.36%  '
$ws.Range('E17').Value = '  -1.94%  '
$ws.Range('D18').Value = '2.594.84'
$ws.Range('E18').Value = '  -3.09%  '
$ws.Range('E19').Value = '  -1.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '346.72'
$ws.Range('E20').Value = '  -2.32%  '
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('E22').Value = '  -2.32%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.33'
$ws.Range('E24').Value = '  -1.60%  '
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('D27').Value = '2.696.40'
$ws.Range('E27').Value = '  -3.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '0.0₃0849'
$ws.Range('E29').Value = '  -3.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.42'
$ws.Range('E30').Value = '  -3.52%  '
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.40'
$ws.Range('E32').Value = '  -2.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '152.94'
$ws.Range('E33').Value = '  -3.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.56'
$ws.Range('E34').Value = '  -2.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.75'
$ws.Range('E35').Value = '  +0.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.05'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.854'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.48'
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('E40').Value = '  -4.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.26'
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '297.95'
$ws.Range('E43').Value = '  -2.47%  '
$ws.Range('E44').Value = '  -3.86%  '
$ws.Range('E45').Value = '  -2.86%  '
$ws.Range('E46').Value = '  -4.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.997'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.82'
$ws.Range('E48').Value = '  -3.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.86'
$ws.Range('E49').Value = '  -3.97%  '
$ws.Range('E50').Value = '  -2.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.29'
$ws.Range('E51').Value = '  +0.27%  '
